$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.532.15'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.820.83'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.37'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5139'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3877'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.80%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08450'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.85'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.111'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.410'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.09'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.517'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.813.08'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001138'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.81'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.79'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.101'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.563.55'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.278'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.13'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.17%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.427'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.024.55'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.78%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.12'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.098'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.29%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.757'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07569'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.674'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2229'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02368'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.205'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.733'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.84%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6340'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.26'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.193'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.54%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.775'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.60%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5938'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '126.01'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.993'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.200'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06989'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.49'
